$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text cells so NDC numbers / dates are not auto-converted by Excel
$textRange = $ws.Range("A2:F3")
$textRange.NumberFormat = "@"

# Row 2: replace old Ascend/Aripiprazole data with Lupin/memantine data (previously on row 3), with updated QTY
$ws.Range("A2").Value = "6818024906"
$ws.Range("B2").Value = "Lupin Pharmaceuticals, Inc."
$ws.Range("C2").Value = "memantine hydrochloride"
$ws.Range("D2").Value = "28 mg/1"
$ws.Range("E2").Value = "H102306"
$ws.Range("F2").Value = "09/23/30"
$ws.Range("G2").Value = 7

# Row 3: new entry for Sun Pharmaceutical / Fenofibrate
$ws.Range("A3").Value = "6330490190"
$ws.Range("B3").Value = "Sun Pharmaceutical Industries, Inc."
$ws.Range("C3").Value = "Fenofibrate"
$ws.Range("D3").Value = "160 mg/1"
$ws.Range("E3").Value = "MHC1672A"
$ws.Range("F3").Value = "11/23/30"
$ws.Range("G3").Value = 9
